$p = $ppt.ActivePresentation

# ---- Slide 19: "Important Driving Factors" ----
$s19 = $p.Slides.Item(19)
$shp19 = $s19.Shapes.Item(4)
$tf19 = $shp19.TextFrame.TextRange

# Paragraph 5 (1-based) - "The columns ..." : remove "/loan_amnt" so that
# "funded_amnt/loan_amnt" becomes "funded_amnt"
$para5 = $tf19.Paragraphs(5, 1)
$para5.Text = 'The columns “funded_amnt”, “int_rate”, “term”, “grade”, “dti”, “verification_status” are decisive factors for identifying the “Charged-off or Fully-paid” loans'

# Paragraph 8 (1-based) - "loan_amnt – Higher the loan amount, likely to be charged-off"
# becomes "Funded_amt – Higher the loan amount, likely to be charged-off"
$para8 = $tf19.Paragraphs(8, 1)
$para8.Text = 'Funded_amt – Higher the loan amount, likely to be charged-off'

# ---- Slide 20: "EDA for the loan dataset revealed that:" ----
$s20 = $p.Slides.Item(20)
$shp20 = $s20.Shapes.Item(4)
$tf20 = $shp20.TextFrame.TextRange

# Paragraph 1 (1-based): merge "hort " + "terms." runs
$para1 = $tf20.Paragraphs(1, 1)
$para1.Text = 'The Bank can reduce the financial loss if it gives loans with short terms.'

# Paragraph 3 (1-based): merge "Average " + "interest rate " + rest
$para3 = $tf20.Paragraphs(3, 1)
$para3.Text = 'Average interest rate of the loans are likely to Fully Paid, so if the banks gives loans at moderate level of interest it can avoid loss.'

# Paragraph 9 (1-based): merge 3 runs into one
$para9 = $tf20.Paragraphs(9, 1)
$para9.Text = 'The Bank can avoid the charged-offs with quality verification process for loan approvals when loan amount is higher for the Grade E, F, G.'
